$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the helper SUM formula in column C (candidate list no longer needs it)
$ws.Range("C1").ClearContents()

# Rewrite the candidate table: drop Glassic/ReStove/Techdalo/Waterfilter,
# fix the casing on Foodgrube/Townbee, and add the new "Operations" candidate.
$ws.Range("A2").Value = "EduGlobe"
$ws.Range("B2").Value = 3

$ws.Range("A3").Value = "Foodgrube"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "Mistub"
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = "SeaSoilution"
$ws.Range("B5").Value = 2

$ws.Range("A6").Value = "Townbee"
$ws.Range("B6").Value = 2

$ws.Range("A7").Value = "Operations"
$ws.Range("B7").Value = 1

# Clear out what used to be rows 8-10
$ws.Range("A8:C10").Clear()

# Highlight the candidate names with a yellow fill
$ws.Range("A2").Interior.Color = 65535
$ws.Range("A4").Interior.Color = 65535
$ws.Range("A5").Interior.Color = 65535
$ws.Range("A6").Interior.Color = 65535
$ws.Range("A7").Interior.Color = 65535

$ws.Range("A7").NumberFormat = "General"
$ws.Range("A3").NumberFormat = "General"

$ws.Range("A7").Select() | Out-Null
